$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.559.22"
$ws.Range("E2").Value = "  +3.71%  "

$ws.Range("D3").Value = "1.603.98"
$ws.Range("E3").Value = "  +3.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.20"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("E6").Value = "  +6.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.96"
$ws.Range("E8").Value = "  +11.03%  "

$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("E10").Value = "  +2.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0599"
$ws.Range("E11").Value = "  +2.74%  "

$ws.Range("E12").Value = "  +2.49%  "

$ws.Range("D13").Value = "1.832.31"
$ws.Range("E13").Value = "  +3.11%  "

$ws.Range("D14").Value = "1.607.07"
$ws.Range("E14").Value = "  +3.29%  "

$ws.Range("D15").Value = "29.607.94"
$ws.Range("E15").Value = "  +3.94%  "

$ws.Range("E16").Value = "  +4.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.75"
$ws.Range("E17").Value = "  +3.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.57"
$ws.Range("E18").Value = "  +3.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.63"
$ws.Range("E19").Value = "  +4.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.63"
$ws.Range("E20").Value = "  +3.57%  "

$ws.Range("D21").Value = "0.0₃0694"
$ws.Range("E21").Value = "  +3.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.02"
$ws.Range("E23").Value = "  +3.36%  "

$ws.Range("E24").Value = "  +3.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.72"
$ws.Range("E26").Value = "  +2.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.31"
$ws.Range("E27").Value = "  +3.64%  "

$ws.Range("E28").Value = "  +5.06%  "

$ws.Range("E29").Value = "  +2.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.06"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.25"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").Value = "1.430.21"
$ws.Range("E34").Value = "  +2.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.11"
$ws.Range("E35").Value = "  +3.63%  "

$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("E37").Value = "  +1.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.83"
$ws.Range("E38").Value = "  +6.91%  "

$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("E41").Value = "  +3.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.96"
$ws.Range("E42").Value = "  +0.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "53.39"
$ws.Range("E43").Value = "  +22.11%  "

$ws.Range("E44").Value = "  +2.29%  "

$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.44"
$ws.Range("E47").Value = "  +1.22%  "

$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("D49").Value = "1.744.00"
$ws.Range("E49").Value = "  +3.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.56"
$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("E51").Value = "  -3.51%  "
